$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header cell format (bold, border, centered) from H1 into the new I1/J1 headers
$ws.Range("H1").Copy($ws.Range("I1:J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# I0 / IF data values for rows 2-76
$values = @(
    @(8,8),
    @(7,7),
    @(8,8),
    @(8,8),
    @(6,6),
    @(8,8),
    @(7,7),
    @(8,8),
    @(8,8),
    @(7,7),
    @(8,8),
    @(7,8),
    @(8,8),
    @(7,7),
    @(7,7),
    @(7,7),
    @(7,8),
    @(8,8),
    @(7,7),
    @(8,8),
    @(10,10),
    @(6,6),
    @(8,8),
    @(8,8),
    @(5,5),
    @(8,8),
    @(7,7),
    @(7,7),
    @(9,9),
    @(7,7),
    @(6,6),
    @(8,8),
    @(7,7),
    @(8,8),
    @(8,9),
    @(8,8),
    @(6,7),
    @(10,10),
    @(10,10),
    @(9,9),
    @(9,9),
    @(8,9),
    @(10,10),
    @(6,6),
    @(7,7),
    @(6,7),
    @(8,8),
    @(8,8),
    @(7,7),
    @(8,8),
    @(8,8),
    @(9,9),
    @(7,7),
    @(8,8),
    @(8,8),
    @(6,6),
    @(8,8),
    @(8,8),
    @(9,9),
    @(8,8),
    @(8,8),
    @(9,9),
    @(9,9),
    @(7,8),
    @(8,9),
    @(8,8),
    @(6,7),
    @(7,7),
    @(6,6),
    @(5,5),
    @(8,8),
    @(9,9),
    @(8,8),
    @(7,7),
    @(5,5)
)

for ($k = 0; $k -lt $values.Length; $k++) {
    $row = $k + 2
    $pair = $values[$k]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
